$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.5901012142502704
$ws.Range("C2").Value = -0.01440167524377749
$ws.Range("D2").Value = 1.061338161608444
$ws.Range("E2").Value = 0.6195339903048283
$ws.Range("F2").Value = -0.1902388420188028
$ws.Range("G2").Value = 0.08989260536700988
$ws.Range("H2").Value = 0.1287596705962219
$ws.Range("B3").Value = 0.1195217185363816
$ws.Range("C3").Value = 1.195261555388603
$ws.Range("D3").Value = 0.7534573840849873
$ws.Range("E3").Value = -0.05631544823864371
$ws.Range("F3").Value = 0.223815999147169
$ws.Range("G3").Value = 0.262683064376381
$ws.Range("B4").Value = 1.270915234191142
$ws.Range("C4").Value = 0.829111062887526
$ws.Range("D4").Value = 0.01933823056389491
$ws.Range("E4").Value = 0.2994696779497076
$ws.Range("F4").Value = 0.3383367431789197
$ws.Range("G4").Value = 0.2296245300636798
$ws.Range("H4").Value = 0.389874912538167
$ws.Range("I4").Value = 0.2978230469237019
$ws.Range("J4").Value = 0.08161590864515353
$ws.Range("B5").Value = 1.510994208950639
$ws.Range("C5").Value = 0.7012213766270079
$ws.Range("D5").Value = 0.9813528240128206
$ws.Range("E5").Value = 1.020219889242033
$ws.Range("F5").Value = 0.9115076761267928
$ws.Range("G5").Value = 1.07175805860128
$ws.Range("H5").Value = 0.9797061929868149
$ws.Range("I5").Value = 0.7634990547082665
$ws.Range("B6").Value = 0.1561095147536435
$ws.Range("C6").Value = 0.4362409621394562
$ws.Range("D6").Value = 0.4751080273686683
$ws.Range("E6").Value = 0.3663958142534284
$ws.Range("F6").Value = 0.5266461967279157
$ws.Range("G6").Value = 0.4345943311134505
$ws.Range("H6").Value = 0.2183871928349021
$ws.Range("B7").Value = 0.4359124617055215
$ws.Range("C7").Value = 0.4747795269347335
$ws.Range("D7").Value = 0.3660673138194937
$ws.Range("E7").Value = 0.5263176962939808
$ws.Range("F7").Value = 0.4342658306795158
$ws.Range("G7").Value = 0.2180586924009674
$ws.Range("B8").Value = 0.588073789631386
$ws.Range("C8").Value = 0.4793615765161462
$ws.Range("D8").Value = 0.6396119589906334
$ws.Range("E8").Value = 0.5475600933761683
$ws.Range("F8").Value = 0.3313529550976199
$ws.Range("G8").Value = 0.01383226865740139
$ws.Range("H8").Value = 0.3060572718689571
$ws.Range("I8").Value = 0.2999120235282561
$ws.Range("B9").Value = 0.1956508791119575
$ws.Range("C9").Value = 0.3559012615864447
$ws.Range("D9").Value = 0.2638493959719796
$ws.Range("E9").Value = 0.04764225769343119
$ws.Range("F9").Value = -0.2698784287467873
$ws.Range("G9").Value = 0.0223465744647684
$ws.Range("H9").Value = 0.01620132612406738
$ws.Range("B10").Value = 0.1455310855428081
$ws.Range("C10").Value = 0.05347921992834301
$ws.Range("D10").Value = -0.1627279183502054
$ws.Range("E10").Value = -0.4802486047904239
$ws.Range("F10").Value = -0.1880236015788682
$ws.Range("G10").Value = -0.1941688499195692
$ws.Range("B11").Value = 0.05627781145257414
$ws.Range("C11").Value = -0.1599293268259743
$ws.Range("D11").Value = -0.4774500132661927
$ws.Range("E11").Value = -0.1852250100546371
$ws.Range("F11").Value = -0.1913702583953381
$ws.Range("B12").Value = -0.2565155703168258
$ws.Range("C12").Value = -0.5740362567570443
$ws.Range("D12").Value = -0.2818112535454886
$ws.Range("E12").Value = -0.2879565018861897
$ws.Range("B13").Value = -0.4995169237785178
$ws.Range("C13").Value = -0.2072919205669621
$ws.Range("D13").Value = -0.2134371689076631
$ws.Range("B14").Value = -0.1459269121925572
$ws.Range("C14").Value = -0.1520721605332582
$ws.Range("B15").Value = -0.1394382194478382
